# Updates cryptos list: price (column D) and 1h volume change (column E)
# values for most rows, plus a swap of the ARBITRUM / InjectiveProtocol
# rows (44 and 45) including their new price/percentage figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = newPrice (optional); E = newPercent }
$updates = @{
    2  = @{ D = '36.200.98';  E = '  -0.76%  ' }
    3  = @{ D = '2.016.06';   E = '  -1.05%  ' }
    4  = @{ E = '  +0.14%  ' }
    5  = @{ D = '251.71';     E = '  +2.99%  ' }
    6  = @{ D = '0.643';      E = '  -2.64%  ' }
    7  = @{ D = '62.92';      E = '  +17.13%  ' }
    8  = @{ E = '  +0.08%  ' }
    9  = @{ D = '59.45';      E = '  -5.39%  ' }
    10 = @{ D = '0.369';      E = '  +1.75%  ' }
    11 = @{ D = '0.0748';     E = '  +0.88%  ' }
    12 = @{ E = '  -0.96%  ' }
    13 = @{ D = '0.935';      E = '  -0.86%  ' }
    14 = @{ D = '14.92';      E = '  +3.42%  ' }
    15 = @{ D = '2.311.68';   E = '  -0.90%  ' }
    16 = @{ D = '5.42';       E = '  +1.78%  ' }
    17 = @{ D = '19.54';      E = '  +14.86%  ' }
    18 = @{ D = '2.026.11';   E = '  -0.59%  ' }
    19 = @{ D = '36.146.79';  E = '  -0.67%  ' }
    20 = @{ D = '72.13';      E = '  +1.61%  ' }
    21 = @{ D = "0.0`u{2083}0856"; E = '  +0.82%  ' }
    22 = @{ D = '5.27';       E = '  +2.46%  ' }
    23 = @{ D = '234.21';     E = '  -1.33%  ' }
    24 = @{ D = '2.75';       E = '  +24.36%  ' }
    25 = @{ E = '  -0.03%  ' }
    26 = @{ E = '  -2.70%  ' }
    27 = @{ D = '9.56';       E = '  +4.71%  ' }
    28 = @{ D = '165.31';     E = '  +1.00%  ' }
    29 = @{ D = '19.62';      E = '  -0.86%  ' }
    30 = @{ D = '5.19';       E = '  +5.04%  ' }
    31 = @{ D = '0.120';      E = '  +0.08%  ' }
    32 = @{ D = '1.20';       E = '  +3.13%  ' }
    33 = @{ D = '0.108';      E = '  +25.40%  ' }
    34 = @{ D = '0.0605';     E = '  +2.09%  ' }
    35 = @{ D = '4.50';       E = '  +2.26%  ' }
    36 = @{ E = '  +11.35%  ' }
    37 = @{ E = '  +0.15%  ' }
    38 = @{ E = '  +0.32%  ' }
    39 = @{ D = '5.75';       E = '  +17.88%  ' }
    40 = @{ D = '0.110';      E = '  +23.24%  ' }
    41 = @{ E = '  -0.37%  ' }
    42 = @{ D = '2.89';       E = '  +1.54%  ' }
    43 = @{ D = '0.0216';     E = '  +1.95%  ' }
}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    if ($vals.ContainsKey('D')) {
        $ws.Cells.Item($rowNum, 4).Value = $vals['D']
    }
    if ($vals.ContainsKey('E')) {
        $ws.Cells.Item($rowNum, 5).Value = $vals['E']
    }
}

# Rows 44 and 45 swap coin identity (ARBITRUM moves to rank 42 / row 44,
# InjectiveProtocol moves to rank 43 / row 45) along with fresh values.
$ws.Cells.Item(44, 2).Value = 'ARBITRUM'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(44, 4).Value = '1.12'
$ws.Cells.Item(44, 5).Value = '  +3.06%  '

$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(45, 4).Value = '16.79'
$ws.Cells.Item(45, 5).Value = '  +6.88%  '
